# Add "Fitting" class data to the MaterialsChart sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns G and H on rows 1 and 2 (material class headers)
$ws.Range("G1").Value = "Carbon Steel"
$ws.Range("H1").Value = "Copper"
$ws.Range("G2").Value = "Copper"
$ws.Range("H2").Value = "Copper"

# Duplicate the existing column F block (rows 1-27) into rows 28-54,
# appending a repeated "Fitting" class section below the original data.
$srcRange = $ws.Range("F1:F27")
$srcRange.Copy()
$destRange = $ws.Range("F28:F54")
$destRange.PasteSpecial(-4104)
